$d = $word.ActiveDocument

# Locate the paragraph that contains the date line "Đà Nẵng, ngày... tháng 10 năm 2025 ..."
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*tháng 10 năm 2025*") {
        $target = $p
        break
    }
}

$pRange = $target.Range.Duplicate

# Find the exact character span of the "10" that must become "12".
$find = $pRange.Duplicate
$find.Find.Execute("10", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$startMonth = $find.Start
$endMonth = $find.End

# Replace the month digits in place.
$monthRange = $d.Range($startMonth, $endMonth)
$monthRange.Text = "12"

# The engine merges adjacent runs that share identical formatting once the
# text is edited, so re-select the surrounding pieces and nudge their
# formatting off/on again to force them back into distinct <w:r> runs that
# match the expected "... tháng " / "12" / " năm 2025 ..." run split.

# Re-locate "ngày" (the text right before " tháng ") so it stays its own run
# instead of merging into the " tháng " run.
$findNgay = $pRange.Duplicate
$findNgay.Find.Execute("ngày", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$ngayStart = $findNgay.Start

# "ngày…" includes the trailing ellipsis character.
$ngayEnd = $ngayStart + 5
$ngayRange = $d.Range($ngayStart, $ngayEnd)
$ngayRange.Font.Bold = 1
$ngayRange.Font.Bold = 0

# Re-apply the same nudge to the "12" run so it stays split from both its
# neighbours (" tháng " before it and " năm 2025 " after it).
$midRange = $d.Range($startMonth, $startMonth + 2)
$midRange.Font.Bold = 1
$midRange.Font.Bold = 0
